$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 1739.5714
$ws.Range("I40").Value = 800
$ws.Range("J40").Value = 2444.25
$ws.Range("K40").Value = 800
$ws.Range("L40").Value = 2444.25
$ws.Range("M40").Value = -625
$ws.Range("N40").Value = -2794.25

# Row 86
$ws.Range("H86").Value = 7900.0713
$ws.Range("I86").Value = 779.9167
$ws.Range("K86").Value = 779.9167
$ws.Range("M86").Value = 343.0833

# Row 89
$ws.Range("H89").Value = 7900.0713
$ws.Range("I89").Value = 779.9167
$ws.Range("K89").Value = 3899.5835
$ws.Range("M89").Value = 1716.4165

# Row 98
$ws.Range("H98").Value = 685
$ws.Range("I98").Value = 685
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 685
$ws.Range("L98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = 813

# Row 122
$ws.Range("H122").Value = 685
$ws.Range("I122").Value = 685
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2055
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = 395

# Row 132
$ws.Range("H132").Value = 5379.647
$ws.Range("I132").Value = 5429.8667
$ws.Range("J132").Value = 5003
$ws.Range("K132").Value = 16289.6001
$ws.Range("L132").Value = 15009
$ws.Range("M132").Value = -13759.6001
$ws.Range("N132").Value = -20069

# Row 133
$ws.Range("H133").Value = 48219.8
$ws.Range("J133").Value = 48219.8
$ws.Range("L133").Value = 48219.8
$ws.Range("N133").Value = -58339.8

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 3457
$ws.Range("I45").Value = 3499.8333
$ws.Range("J45").Value = 3420.2856
$ws.Range("K45").Value = 3499.8333
$ws.Range("L45").Value = 3420.2856
$ws.Range("M45").Value = -3122.8333
$ws.Range("N45").Value = -4174.2856

# Row 61
$ws.Range("H61").Value = 3583.7827
$ws.Range("I61").Value = 3520.3333
$ws.Range("K61").Value = 3520.3333
$ws.Range("M61").Value = -3308.3333

# Row 132
$ws.Range("H132").Value = 29590.79
$ws.Range("I132").Value = 3161.0833
$ws.Range("K132").Value = 9483.249899999999
$ws.Range("M132").Value = -6953.249899999999

# Row 136
$ws.Range("H136").Value = 3583.7827
$ws.Range("I136").Value = 3520.3333
$ws.Range("K136").Value = 10560.9999
$ws.Range("M136").Value = -8010.999899999999

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 3837.2856
$ws.Range("I134").Value = 3942.3704
$ws.Range("K134").Value = 11827.1112
$ws.Range("M134").Value = -9292.111199999999

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2780.4614
$ws.Range("I31").Value = 1153.909
$ws.Range("J31").Value = 3419.4644
$ws.Range("K31").Value = 1153.909
$ws.Range("L31").Value = 3419.4644
$ws.Range("M31").Value = -858.9090000000001
$ws.Range("N31").Value = -4009.4644

# Row 34
$ws.Range("H34").Value = 2780.4614
$ws.Range("I34").Value = 1153.909
$ws.Range("J34").Value = 3419.4644
$ws.Range("K34").Value = 1153.909
$ws.Range("L34").Value = 3419.4644
$ws.Range("M34").Value = -951.9090000000001
$ws.Range("N34").Value = -3823.4644

# Row 132
$ws.Range("H132").Value = 4884.6665
$ws.Range("I132").Value = 2025
$ws.Range("K132").Value = 6075
$ws.Range("M132").Value = -3545

# Row 134
$ws.Range("H134").Value = 1312.7727
$ws.Range("I134").Value = 1247.625
$ws.Range("J134").Value = 1486.5
$ws.Range("K134").Value = 3742.875
$ws.Range("L134").Value = 4459.5
$ws.Range("M134").Value = -1207.875
$ws.Range("N134").Value = -9529.5

$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Range("H122").Value = 508.7647
$ws.Range("I122").Value = 362.6
$ws.Range("J122").Value = 569.6667
$ws.Range("K122").Value = 3263.4
$ws.Range("L122").Value = 5127.0003
$ws.Range("M122").Value = -813.4000000000001
$ws.Range("N122").Value = -10027.0003

# Row 123
$ws.Range("H123").Value = 4499.8
$ws.Range("I123").Value = 1029
$ws.Range("J123").Value = 5367.5
$ws.Range("K123").Value = 3087
$ws.Range("L123").Value = 16102.5
$ws.Range("M123").Value = -637
$ws.Range("N123").Value = -21002.5

# Row 131
$ws.Range("H131").Value = 711.36
$ws.Range("I131").Value = 498
$ws.Range("J131").Value = 722.5895
$ws.Range("K131").Value = 1494
$ws.Range("L131").Value = 2167.7685
$ws.Range("M131").Value = 3546
$ws.Range("N131").Value = -12247.7685

# Row 139
$ws.Range("H139").Value = 2022.92
$ws.Range("I139").Value = 1353.5264
$ws.Range("J139").Value = 4142.6665
$ws.Range("K139").Value = 4060.5792
$ws.Range("L139").Value = 12427.9995
$ws.Range("M139").Value = 1079.4208
$ws.Range("N139").Value = -22707.9995

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3756.5217
$ws.Range("I80").Value = 2800
$ws.Range("J80").Value = 4492.3076
$ws.Range("K80").Value = 2800
$ws.Range("L80").Value = 4492.3076
$ws.Range("M80").Value = -1802
$ws.Range("N80").Value = -6488.3076

# Row 83
$ws.Range("H83").Value = 3756.5217
$ws.Range("I83").Value = 2800
$ws.Range("J83").Value = 4492.3076
$ws.Range("K83").Value = 14000
$ws.Range("L83").Value = 22461.538
$ws.Range("M83").Value = -9008
$ws.Range("N83").Value = -32445.538

# Row 113
$ws.Range("H113").Value = 2194.5908
$ws.Range("I113").Value = 1684.5333
$ws.Range("J113").Value = 3287.5715
$ws.Range("K113").Value = 1684.5333
$ws.Range("L113").Value = 3287.5715
$ws.Range("M113").Value = 485.4666999999999
$ws.Range("N113").Value = -7627.5715

# Row 132
$ws.Range("H132").Value = 23796.652
$ws.Range("I132").Value = 1460.3636
$ws.Range("J132").Value = 44271.582
$ws.Range("K132").Value = 4381.0908
$ws.Range("L132").Value = 132814.746
$ws.Range("M132").Value = -1851.0908
$ws.Range("N132").Value = -137874.746

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 4541.6
$ws.Range("I61").Value = 1728.125
$ws.Range("K61").Value = 1728.125
$ws.Range("M61").Value = -1526.125

# Row 68
$ws.Range("H68").Value = 2419.4
$ws.Range("J68").Value = 2999
$ws.Range("L68").Value = 2999
$ws.Range("N68").Value = -4497

# Row 71
$ws.Range("H71").Value = 2419.4
$ws.Range("J71").Value = 2999
$ws.Range("L71").Value = 14995
$ws.Range("N71").Value = -22483

# Row 100
$ws.Range("H100").Value = 2366.3333
$ws.Range("I100").Value = 1950
$ws.Range("J100").Value = 2574.5
$ws.Range("K100").Value = 1950
$ws.Range("L100").Value = 2574.5
$ws.Range("M100").Value = -1409
$ws.Range("N100").Value = -3656.5

# Row 113
$ws.Range("H113").Value = 4541.6
$ws.Range("I113").Value = 1728.125
$ws.Range("K113").Value = 1728.125
$ws.Range("M113").Value = 441.875

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 1186.9166
$ws.Range("I81").Value = 1171.5555
$ws.Range("J81").Value = 1233
$ws.Range("K81").Value = 2343.111
$ws.Range("L81").Value = 2466
$ws.Range("M81").Value = -1282.111
$ws.Range("N81").Value = -4588

# Row 84
$ws.Range("H84").Value = 1186.9166
$ws.Range("I84").Value = 1171.5555
$ws.Range("J84").Value = 1233
$ws.Range("K84").Value = 11715.555
$ws.Range("L84").Value = 12330
$ws.Range("M84").Value = -6411.555
$ws.Range("N84").Value = -22938

# Row 126
$ws.Range("H126").Value = 1543.0714
$ws.Range("I126").Value = 1569.4615
$ws.Range("J126").Value = 1200
$ws.Range("K126").Value = 4708.3845
$ws.Range("L126").Value = 3600
$ws.Range("M126").Value = -2238.3845
$ws.Range("N126").Value = -8540

# Row 132
$ws.Range("H132").Value = 1165.7667
$ws.Range("I132").Value = 715.9474
$ws.Range("J132").Value = 1942.7273
$ws.Range("K132").Value = 2147.8422
$ws.Range("L132").Value = 5828.1819
$ws.Range("M132").Value = 382.1578
$ws.Range("N132").Value = -10888.1819
